# Update "Forecast Comparison" sheet with corrected forecast output:
#  - insert a new "Week_Start_Date" column after "Week" (shifts ASIN.. right by one)
#  - rewrite the Week labels without the leading zero (W01 -> W1, ... W09 -> W9)
#  - populate the new Week_Start_Date column with the week's start date (as text)
#  - correct a few MyForecast values
#  - store is_holiday_week as a real boolean instead of 0/1 numbers

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column before the current column B (ASIN), shifting
#    ASIN / MyForecast / Amazon*Forecast / Product Title / is_holiday_week right by one.
$ws.Columns("B").Insert()

# 2) New column header
$ws.Range("B1").Value = "Week_Start_Date"

# 3) Per-week data: row, new Week label, Week_Start_Date text
$weekInfo = @(
    @{ Row = 2;  Week = "W1";  Start = "2025-01-05" },
    @{ Row = 3;  Week = "W2";  Start = "2025-01-12" },
    @{ Row = 4;  Week = "W3";  Start = "2025-01-19" },
    @{ Row = 5;  Week = "W4";  Start = "2025-01-26" },
    @{ Row = 6;  Week = "W5";  Start = "2025-02-02" },
    @{ Row = 7;  Week = "W6";  Start = "2025-02-09" },
    @{ Row = 8;  Week = "W7";  Start = "2025-02-16" },
    @{ Row = 9;  Week = "W8";  Start = "2025-02-23" },
    @{ Row = 10; Week = "W9";  Start = "2025-03-02" },
    @{ Row = 11; Week = "W10"; Start = "2025-03-09" },
    @{ Row = 12; Week = "W11"; Start = "2025-03-16" },
    @{ Row = 13; Week = "W12"; Start = "2025-03-23" },
    @{ Row = 14; Week = "W13"; Start = "2025-03-30" },
    @{ Row = 15; Week = "W14"; Start = "2025-04-06" },
    @{ Row = 16; Week = "W15"; Start = "2025-04-13" },
    @{ Row = 17; Week = "W16"; Start = "2025-04-20" }
)

foreach ($info in $weekInfo) {
    $r = $info.Row

    # Week label (A) without leading zero
    $ws.Cells.Item($r, 1).Value = $info.Week

    # Week_Start_Date (B) - force text so Excel doesn't coerce to a date serial
    $cell = $ws.Cells.Item($r, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $info.Start

    # is_holiday_week (J) - store as real boolean FALSE
    $ws.Cells.Item($r, 10).Value = $false
}

# 4) Corrected MyForecast values (column D after the insert)
$ws.Cells.Item(2, 4).Value = 21
$ws.Cells.Item(14, 4).Value = 19
$ws.Cells.Item(15, 4).Value = 17
